# ------------------------------------------------------------------
# Adds a new "2022-Q4" quarterly sheet to the workbook (inserted right
# after "总计" and before "2022-Q3"), fills it with the fund-holding
# table for that quarter, and updates the "总计" summary sheet with
# the corresponding aggregate row (shifting the existing 2022-Q3 /
# 2022-Q2 rows down).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Insert the new "2022-Q4" worksheet before "2022-Q3" --------
$refSheet = $wb.Worksheets.Item(2)
$new = $wb.Worksheets.Add($refSheet)
$new.Name = "2022-Q4"

# Header row (same layout/style as the other quarterly sheets)
$new.Cells.Item(1,2).Value = "基金代码"
$new.Cells.Item(1,3).Value = "基金名称"
$new.Cells.Item(1,4).Value = "基金规模"
$new.Cells.Item(1,5).Value = "股票总仓位"
$new.Cells.Item(1,6).Value = "仓位占比"
$new.Cells.Item(1,7).Value = "持有市值(亿元)"
$new.Cells.Item(1,8).Value = "仓位排名"

# ---- 2. Fill the quarterly fund-holding data ------------------------
$new.Cells.Item(2,1).Value = 0
$new.Cells.Item(2,2).Value = "'012382"
$new.Cells.Item(2,3).Value = "泰达宏利新兴景气龙头混合A"
$new.Cells.Item(2,4).Value = "'16.52"
$new.Cells.Item(2,5).Value = "'90.73"
$new.Cells.Item(2,6).Value = "'4.41"
$new.Cells.Item(2,7).Value = "'0.7285"
$new.Cells.Item(2,8).Value = 10
$new.Cells.Item(3,1).Value = 1
$new.Cells.Item(3,2).Value = "'501081"
$new.Cells.Item(3,3).Value = "中欧科创主题混合（LOF）A"
$new.Cells.Item(3,4).Value = "'7.06"
$new.Cells.Item(3,5).Value = "'87.64"
$new.Cells.Item(3,6).Value = "'2.82"
$new.Cells.Item(3,7).Value = "'0.1991"
$new.Cells.Item(3,8).Value = 10
$new.Cells.Item(4,1).Value = 2
$new.Cells.Item(4,2).Value = "'011506"
$new.Cells.Item(4,3).Value = "建信高端装备股票A"
$new.Cells.Item(4,4).Value = "'3.84"
$new.Cells.Item(4,5).Value = "'89.26"
$new.Cells.Item(4,6).Value = "'3.42"
$new.Cells.Item(4,7).Value = "'0.1313"
$new.Cells.Item(4,8).Value = 3
$new.Cells.Item(5,1).Value = 3
$new.Cells.Item(5,2).Value = "'011793"
$new.Cells.Item(5,3).Value = "建信智能汽车股票"
$new.Cells.Item(5,4).Value = "'4.88"
$new.Cells.Item(5,5).Value = "'87.61"
$new.Cells.Item(5,6).Value = "'2.17"
$new.Cells.Item(5,7).Value = "'0.1059"
$new.Cells.Item(5,8).Value = 9
$new.Cells.Item(6,1).Value = 4
$new.Cells.Item(6,2).Value = "'012383"
$new.Cells.Item(6,3).Value = "泰达宏利新兴景气龙头混合C"
$new.Cells.Item(6,4).Value = "'2.19"
$new.Cells.Item(6,5).Value = "'90.73"
$new.Cells.Item(6,6).Value = "'4.41"
$new.Cells.Item(6,7).Value = "'0.0966"
$new.Cells.Item(6,8).Value = 10
$new.Cells.Item(7,1).Value = 5
$new.Cells.Item(7,2).Value = "'011093"
$new.Cells.Item(7,3).Value = "永赢宏泽一年定期开放灵活配置混合"
$new.Cells.Item(7,4).Value = "'14.98"
$new.Cells.Item(7,5).Value = "'48.20"
$new.Cells.Item(7,6).Value = "'0.59"
$new.Cells.Item(7,7).Value = "'0.0884"
$new.Cells.Item(7,8).Value = 4
$new.Cells.Item(8,1).Value = 6
$new.Cells.Item(8,2).Value = "'011507"
$new.Cells.Item(8,3).Value = "建信高端装备股票C"
$new.Cells.Item(8,4).Value = "'0.90"
$new.Cells.Item(8,5).Value = "'89.26"
$new.Cells.Item(8,6).Value = "'3.42"
$new.Cells.Item(8,7).Value = "'0.0308"
$new.Cells.Item(8,8).Value = 3
$new.Cells.Item(9,1).Value = 7
$new.Cells.Item(9,2).Value = "'011351"
$new.Cells.Item(9,3).Value = "金鹰年年邮益一年持有期混合A"
$new.Cells.Item(9,4).Value = "'3.04"
$new.Cells.Item(9,5).Value = "'39.17"
$new.Cells.Item(9,6).Value = "'0.75"
$new.Cells.Item(9,7).Value = "'0.0228"
$new.Cells.Item(9,8).Value = 9
$new.Cells.Item(10,1).Value = 8
$new.Cells.Item(10,2).Value = "'014478"
$new.Cells.Item(10,3).Value = "中加低碳经济六个月持有期混合A"
$new.Cells.Item(10,4).Value = "'0.61"
$new.Cells.Item(10,5).Value = "'93.20"
$new.Cells.Item(10,6).Value = "'3.63"
$new.Cells.Item(10,7).Value = "'0.0221"
$new.Cells.Item(10,8).Value = 7
$new.Cells.Item(11,1).Value = 9
$new.Cells.Item(11,2).Value = "'006836"
$new.Cells.Item(11,3).Value = "永赢惠泽一年定期开放灵活配置混合"
$new.Cells.Item(11,4).Value = "'3.73"
$new.Cells.Item(11,5).Value = "'48.14"
$new.Cells.Item(11,6).Value = "'0.59"
$new.Cells.Item(11,7).Value = "'0.0220"
$new.Cells.Item(11,8).Value = 4
$new.Cells.Item(12,1).Value = 10
$new.Cells.Item(12,2).Value = "'008061"
$new.Cells.Item(12,3).Value = "惠升惠新灵活配置混合A"
$new.Cells.Item(12,4).Value = "'0.37"
$new.Cells.Item(12,5).Value = "'89.09"
$new.Cells.Item(12,6).Value = "'4.87"
$new.Cells.Item(12,7).Value = "'0.0180"
$new.Cells.Item(12,8).Value = 5
$new.Cells.Item(13,1).Value = 11
$new.Cells.Item(13,2).Value = "'014479"
$new.Cells.Item(13,3).Value = "中加低碳经济六个月持有期混合C"
$new.Cells.Item(13,4).Value = "'0.17"
$new.Cells.Item(13,5).Value = "'93.20"
$new.Cells.Item(13,6).Value = "'3.63"
$new.Cells.Item(13,7).Value = "'0.0062"
$new.Cells.Item(13,8).Value = 7
$new.Cells.Item(14,1).Value = 12
$new.Cells.Item(14,2).Value = "'970046"
$new.Cells.Item(14,3).Value = "东海证券海睿健行灵活配置混合A"
$new.Cells.Item(14,4).Value = "'0.14"
$new.Cells.Item(14,5).Value = "'87.16"
$new.Cells.Item(14,6).Value = "'3.34"
$new.Cells.Item(14,7).Value = "'0.0047"
$new.Cells.Item(14,8).Value = 10
$new.Cells.Item(15,1).Value = 13
$new.Cells.Item(15,2).Value = "'970083"
$new.Cells.Item(15,3).Value = "东海证券海盈6个月持有期混合"
$new.Cells.Item(15,4).Value = "'0.09"
$new.Cells.Item(15,5).Value = "'37.93"
$new.Cells.Item(15,6).Value = "'3.35"
$new.Cells.Item(15,7).Value = "'0.0030"
$new.Cells.Item(15,8).Value = 3
$new.Cells.Item(16,1).Value = 14
$new.Cells.Item(16,2).Value = "'970047"
$new.Cells.Item(16,3).Value = "东海证券海睿健行灵活配置混合B"
$new.Cells.Item(16,4).Value = "'0.09"
$new.Cells.Item(16,5).Value = "'87.16"
$new.Cells.Item(16,6).Value = "'3.34"
$new.Cells.Item(16,7).Value = "'0.0030"
$new.Cells.Item(16,8).Value = 10
$new.Cells.Item(17,1).Value = 15
$new.Cells.Item(17,2).Value = "'008062"
$new.Cells.Item(17,3).Value = "惠升惠新灵活配置混合C"
$new.Cells.Item(17,4).Value = "'0.06"
$new.Cells.Item(17,5).Value = "'89.09"
$new.Cells.Item(17,6).Value = "'4.87"
$new.Cells.Item(17,7).Value = "'0.0029"
$new.Cells.Item(17,8).Value = 5
$new.Cells.Item(18,1).Value = 16
$new.Cells.Item(18,2).Value = "'007533"
$new.Cells.Item(18,3).Value = "格林创新成长混合A"
$new.Cells.Item(18,4).Value = "'0.05"
$new.Cells.Item(18,5).Value = "'75.81"
$new.Cells.Item(18,6).Value = "'4.55"
$new.Cells.Item(18,7).Value = "'0.0023"
$new.Cells.Item(18,8).Value = 8
$new.Cells.Item(19,1).Value = 17
$new.Cells.Item(19,2).Value = "'007534"
$new.Cells.Item(19,3).Value = "格林创新成长混合C"
$new.Cells.Item(19,4).Value = "'0.05"
$new.Cells.Item(19,5).Value = "'75.81"
$new.Cells.Item(19,6).Value = "'4.55"
$new.Cells.Item(19,7).Value = "'0.0023"
$new.Cells.Item(19,8).Value = 8
$new.Cells.Item(20,1).Value = 18
$new.Cells.Item(20,2).Value = "'011352"
$new.Cells.Item(20,3).Value = "金鹰年年邮益一年持有期混合C"
$new.Cells.Item(20,4).Value = "'0.23"
$new.Cells.Item(20,5).Value = "'39.17"
$new.Cells.Item(20,6).Value = "'0.75"
$new.Cells.Item(20,7).Value = "'0.0017"
$new.Cells.Item(20,8).Value = 9
$new.Cells.Item(21,1).Value = 19
$new.Cells.Item(21,2).Value = "'008533"
$new.Cells.Item(21,3).Value = "惠升惠兴混合A"
$new.Cells.Item(21,4).Value = "'0.03"
$new.Cells.Item(21,5).Value = "'23.29"
$new.Cells.Item(21,6).Value = "'1.93"
$new.Cells.Item(21,7).Value = "'0.0006"
$new.Cells.Item(21,8).Value = 6
$new.Cells.Item(22,1).Value = 20
$new.Cells.Item(22,2).Value = "'017290"
$new.Cells.Item(22,3).Value = "中欧科创主题混合（LOF）C"
$new.Cells.Item(22,4).Value = "'0.00"
$new.Cells.Item(22,5).Value = "'87.64"
$new.Cells.Item(22,6).Value = "'2.82"
$new.Cells.Item(22,7).Value = 0
$new.Cells.Item(22,8).Value = 10
$new.Cells.Item(23,1).Value = 21
$new.Cells.Item(23,2).Value = "'008534"
$new.Cells.Item(23,3).Value = "惠升惠兴混合C"
$new.Cells.Item(23,4).Value = "'0.00"
$new.Cells.Item(23,5).Value = "'23.29"
$new.Cells.Item(23,6).Value = "'1.93"
$new.Cells.Item(23,7).Value = 0
$new.Cells.Item(23,8).Value = 6

# ---- 3. Match styling used by the other quarterly sheets ------------
# Header row style (bold + border), copied from the "总计" header cell.
$headerStyleSrc = $wb.Worksheets.Item(1).Range("B1")
$headerStyleSrc.Copy()
$new.Range("B1:H1").PasteSpecial(-4122)

# Row-index column (A) style, copied from the "总计" index column.
$indexStyleSrc = $wb.Worksheets.Item(1).Range("A2")
$indexStyleSrc.Copy()
$new.Range("A2:A23").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- 4. Update the "总计" summary sheet ------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Shift the existing "2022-Q2" row (was row 3) down to row 4.
$ws1.Cells.Item(4,1).Value = 2
$ws1.Cells.Item(4,2).Value = "2022-Q2"
$ws1.Cells.Item(4,3).Value = 8
$ws1.Cells.Item(4,4).Value = 1.2

# Shift the existing "2022-Q3" row (was row 2) down to row 3.
$ws1.Cells.Item(3,1).Value = 1
$ws1.Cells.Item(3,2).Value = "2022-Q3"
$ws1.Cells.Item(3,3).Value = 6
$ws1.Cells.Item(3,4).Value = 0.12

# New "2022-Q4" aggregate row goes into row 2.
$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,2).Value = "2022-Q4"
$ws1.Cells.Item(2,3).Value = 22
$ws1.Cells.Item(2,4).Value = 1.49

# Make sure the new row-4 index cell (A4) carries the same style as
# the other index cells in column A.
$indexStyleSrc2 = $ws1.Range("A2")
$indexStyleSrc2.Copy()
$ws1.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- 5. Restore the originally-active "2022-Q2" tab ------------------
$wb.Worksheets.Item(4).Activate()
